$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.591.03'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '2.286.37'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '313.36'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").Value = '105.05'
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("E7").Value = '  -0.62%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  -0.75%  '
$ws.Range("D10").Value = '39.58'
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").Value = '0.0903'
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = '8.38'
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D14").Value = '0.992'
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("D15").Value = '15.20'
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '2.632.15'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").Value = '2.276.97'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '42.769.61'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = '7.37'
$ws.Range("E19").Value = '  -1.25%  '
$ws.Range("D20").Value = '13.80'
$ws.Range("E20").Value = '  +22.28%  '
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").Value = '73.79'
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").Value = '264.92'
$ws.Range("E24").Value = '  -4.19%  '
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  -2.88%  '
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").Value = '10.84'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").Value = '7.19'
$ws.Range("E28").Value = '  +22.32%  '
$ws.Range("D29").Value = '2.35'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("D31").Value = '37.06'
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("D32").Value = '166.96'
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("D33").Value = '0.0871'
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").Value = '0.131'
$ws.Range("E34").Value = '  -3.10%  '
$ws.Range("D35").Value = '2.58'
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("E36").Value = '  -4.36%  '
$ws.Range("D37").Value = '4.53'
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("E38").Value = '  -4.27%  '
$ws.Range("D39").Value = '3.78'
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("D40").Value = '2.66'
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("E41").Value = '  +4.77%  '
$ws.Range("D42").Value = '70.37'
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("E43").Value = '  +1.69%  '
$ws.Range("D44").Value = '94.72'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").Value = '12.15'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '113.45'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.726.38'
$ws.Range("E48").Value = '  +8.57%  '
$ws.Range("D49").Value = '79.16'
$ws.Range("E49").Value = '  -3.16%  '
$ws.Range("D50").Value = '8.73'
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("E51").Value = '  -0.61%  '
